$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.746.08"
$ws.Range("E2").Value = "  +0.51%  "

# Row 3
$ws.Range("D3").Value = "'2.644.96"
$ws.Range("E3").Value = "  -1.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").Value = "'608.31"
$ws.Range("E5").Value = "  -0.16%  "

# Row 6
$ws.Range("D6").Value = "'147.53"
$ws.Range("E6").Value = "  +2.61%  "

# Row 7
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("E8").Value = "  +0.18%  "

# Row 9
$ws.Range("E9").Value = "  +1.98%  "

# Row 10
$ws.Range("D10").Value = "'0.384"
$ws.Range("E10").Value = "  +6.92%  "

# Row 11
$ws.Range("E11").Value = "  -1.14%  "

# Row 12
$ws.Range("D12").Value = "'0.151"
$ws.Range("E12").Value = "  -0.87%  "

# Row 13
$ws.Range("D13").Value = "'27.49"
$ws.Range("E13").Value = "  +0.11%  "

# Row 14
$ws.Range("D14").Value = "'3.123.25"
$ws.Range("E14").Value = "  -0.75%  "

# Row 15
$ws.Range("D15").Value = "'63.648.87"
$ws.Range("E15").Value = "  +0.60%  "

# Row 16
$ws.Range("E16").Value = "  +1.03%  "

# Row 17
$ws.Range("D17").Value = "'2.663.93"
$ws.Range("E17").Value = "  +0.82%  "

# Row 18
$ws.Range("D18").Value = "'11.79"
$ws.Range("E18").Value = "  +2.66%  "

# Row 19
$ws.Range("D19").Value = "'4.57"
$ws.Range("E19").Value = "  +3.79%  "

# Row 20
$ws.Range("D20").Value = "'346.59"
$ws.Range("E20").Value = "  +1.68%  "

# Row 21
$ws.Range("E21").Value = "  +0.33%  "

# Row 22
$ws.Range("E22").Value = "  -0.16%  "

# Row 23
$ws.Range("E23").Value = "  -2.17%  "

# Row 24
$ws.Range("D24").Value = "'66.35"
$ws.Range("E24").Value = "  -1.84%  "

# Row 25
$ws.Range("E25").Value = "  +6.66%  "

# Row 26
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'9.22"
$ws.Range("E26").Value = "  +7.26%  "

# Row 27
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").Value = "'1.68"
$ws.Range("E27").Value = "  +1.71%  "

# Row 28
$ws.Range("D28").Value = "'560.13"
$ws.Range("E28").Value = "  +4.23%  "

# Row 29
$ws.Range("D29").Value = "'8.11"
$ws.Range("E29").Value = "  +3.01%  "

# Row 30
$ws.Range("D30").Value = "'0.163"
$ws.Range("E30").Value = "  -1.43%  "

# Row 31
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.08%  "

# Row 32
$ws.Range("D32").Value = "'2.04"
$ws.Range("E32").Value = "  -0.73%  "

# Row 33
$ws.Range("D33").Value = "'0.0₃0850"
$ws.Range("E33").Value = "  +4.95%  "

# Row 34
$ws.Range("E34").Value = "  -2.33%  "

# Row 35
$ws.Range("D35").Value = "'5.29"
$ws.Range("E35").Value = "  +3.85%  "

# Row 36
$ws.Range("D36").Value = "'168.89"
$ws.Range("E36").Value = "  -2.19%  "

# Row 37
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.06%  "

# Row 38
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.405"
$ws.Range("E38").Value = "  -0.21%  "

# Row 39
$ws.Range("E39").Value = "  +5.14%  "

# Row 40
$ws.Range("E40").Value = "  -0.46%  "

# Row 41
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.09%  "

# Row 42
$ws.Range("D42").Value = "'165.19"
$ws.Range("E42").Value = "  -5.48%  "

# Row 43
$ws.Range("D43").Value = "'40.05"
$ws.Range("E43").Value = "  -0.39%  "

# Row 44
$ws.Range("D44").Value = "'3.81"
$ws.Range("E44").Value = "  +1.51%  "

# Row 45
$ws.Range("D45").Value = "'21.95"
$ws.Range("E45").Value = "  -1.03%  "

# Row 46
$ws.Range("D46").Value = "'0.0566"
$ws.Range("E46").Value = "  +0.13%  "

# Row 47
$ws.Range("D47").Value = "'0.629"
$ws.Range("E47").Value = "  -1.03%  "

# Row 48
$ws.Range("D48").Value = "'2.02"
$ws.Range("E48").Value = "  +15.64%  "

# Row 49
$ws.Range("D49").Value = "'0.0244"
$ws.Range("E49").Value = "  +1.50%  "

# Row 50
$ws.Range("D50").Value = "'0.0959"
$ws.Range("E50").Value = "  -0.64%  "

# Row 51
$ws.Range("D51").Value = "'18.85"
$ws.Range("E51").Value = "  -0.34%  "
